$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "68.511.85", "0.999") are preserved exactly as text, matching
# the source data which stores these as literal strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.511.85"
$ws.Range("E2").Value = "  -2.57%  "
$ws.Range("D3").Value = "3.472.37"
$ws.Range("E3").Value = "  -3.99%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "578.52"
$ws.Range("E5").Value = "  -3.94%  "
$ws.Range("D6").Value = "189.03"
$ws.Range("E6").Value = "  -3.48%  "
$ws.Range("D7").Value = "3.458.20"
$ws.Range("E7").Value = "  -4.06%  "
$ws.Range("D8").Value = "0.599"
$ws.Range("E8").Value = "  -4.45%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "0.200"
$ws.Range("E10").Value = "  -5.83%  "
$ws.Range("D11").Value = "0.610"
$ws.Range("E11").Value = "  -5.60%  "
$ws.Range("D12").Value = "51.35"
$ws.Range("E12").Value = "  -3.63%  "
$ws.Range("D13").Value = "0.0000282"
$ws.Range("E13").Value = "  -7.48%  "
$ws.Range("D14").Value = "9.04"
$ws.Range("E14").Value = "  -5.46%  "
$ws.Range("D15").Value = "4.013.67"
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("D16").Value = "633.24"
$ws.Range("E16").Value = "  +4.97%  "
$ws.Range("D17").Value = "68.775.76"
$ws.Range("E17").Value = "  -2.34%  "
$ws.Range("D18").Value = "3.478.52"
$ws.Range("E18").Value = "  -4.27%  "
$ws.Range("D19").Value = "12.31"
$ws.Range("E19").Value = "  -4.92%  "
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("E21").Value = "  -5.53%  "
$ws.Range("D22").Value = "0.940"
$ws.Range("E22").Value = "  -6.02%  "
$ws.Range("D23").Value = "17.61"
$ws.Range("E23").Value = "  -4.97%  "
$ws.Range("D24").Value = "5.36"
$ws.Range("E24").Value = "  +2.87%  "
$ws.Range("D25").Value = "99.08"
$ws.Range("E25").Value = "  -4.33%  "
$ws.Range("E26").Value = "  -7.28%  "
$ws.Range("D27").Value = "2.85"
$ws.Range("E27").Value = "  -4.95%  "
$ws.Range("E28").Value = "  +2.60%  "
$ws.Range("D29").Value = "9.98"
$ws.Range("E29").Value = "  -5.61%  "
$ws.Range("D30").Value = "9.13"
$ws.Range("E30").Value = "  -5.57%  "
$ws.Range("D31").Value = "32.32"
$ws.Range("E31").Value = "  -4.48%  "
$ws.Range("D32").Value = "6.66"
$ws.Range("E32").Value = "  -8.41%  "
$ws.Range("D33").Value = "4.05"
$ws.Range("E33").Value = "  -13.45%  "
$ws.Range("D34").Value = "11.53"
$ws.Range("E34").Value = "  -6.19%  "
$ws.Range("E35").Value = "  -7.60%  "
$ws.Range("D36").Value = "60.83"
$ws.Range("E36").Value = "  -3.77%  "
$ws.Range("D37").Value = "3.684.53"
$ws.Range("E37").Value = "  -6.32%  "
$ws.Range("D38").Value = "0.996"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").Value = "0.0₃0789"
$ws.Range("E39").Value = "  -10.49%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "3.55"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "500.05"
$ws.Range("E41").Value = "  -5.96%  "
$ws.Range("D42").Value = "2.90"
$ws.Range("E42").Value = "  -5.17%  "
$ws.Range("D43").Value = "0.364"
$ws.Range("E43").Value = "  -6.48%  "
$ws.Range("D44").Value = "0.132"
$ws.Range("E44").Value = "  -1.72%  "
$ws.Range("D45").Value = "34.03"
$ws.Range("E45").Value = "  -7.73%  "
$ws.Range("D46").Value = "0.0437"
$ws.Range("E46").Value = "  -5.39%  "
$ws.Range("D47").Value = "3.30"
$ws.Range("E47").Value = "  -8.41%  "
$ws.Range("E48").Value = "  -3.52%  "
$ws.Range("E49").Value = "  -4.95%  "
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").Value = "8.05"
$ws.Range("E51").Value = "  -6.17%  "
